$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, shifting existing rows 217:312 down to 218:313
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly record
$row = 217
$ws.Cells.Item($row, 1).Value = 6
$ws.Cells.Item($row, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item($row, 3).Value = 'Metropolitana'
$ws.Cells.Item($row, 4).Value = 44825
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = 100112026
$ws.Cells.Item($row, 7).Value = 'Haba'
$ws.Cells.Item($row, 8).Value = 'Sin especificar'
$ws.Cells.Item($row, 9).Value = 'Primera'
$ws.Cells.Item($row, 10).Value = 650
$ws.Cells.Item($row, 11).Value = 9000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 9585
$ws.Cells.Item($row, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($row, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item($row, 16).Value = 383
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = 'Hortaliza'
